$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Page Objects And Actions")
$ws.Activate()

# Column D on this sheet (DriverViewPage action list) gets a new entry
# "refreshPage" inserted at D30; every existing entry from D30 down to
# D162 shifts down one row (to D31:D163). Columns A-C are untouched.
# Shift from the bottom up so we never clobber a value before reading it.
for ($r = 163; $r -ge 31; $r--) {
    $srcText = $ws.Cells.Item($r - 1, 4).Text
    $ws.Cells.Item($r, 4).Value = $srcText
}
$ws.Cells.Item(30, 4).Value = "refreshPage"

# Re-apply the existing sort on the (now one-row-longer) D2:D163 range so
# the sheet's recorded sortState/sortCondition track the new extent.
$sort = $ws.Sort
$sort.SortFields.Clear()
$sort.SortFields.Add($ws.Range("D2:D163"))
$sort.SetRange($ws.Range("D2:D163"))
$sort.Header = 2
$sort.Apply()

# Update the view: scroll position/top-left cell and the active selection,
# which now spans the whole (longer) D column data range.
$win = $excel.ActiveWindow
$win.ScrollRow = 148
$win.ScrollColumn = 1
$ws.Range("D2:D163").Select()
